$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-9
# from 2023-10-08 (45207) to 2023-10-09 (45208), keeping existing formatting.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45208
}
